$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 336, pushing the existing rows 336-351 down to 337-352.
$ws.Rows.Item(336).Insert()

# Populate the newly inserted row 336 with a new weekly data point for
# Berenjena / Vega Central Mapocho de Santiago (same dimensions as the row
# that used to sit at 336, now shifted to 337, but with its own
# date/volume/average-price values).
$ws.Cells.Item(336, 1).Value = 9
$ws.Cells.Item(336, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(336, 3).Value = "Metropolitana"
$ws.Cells.Item(336, 4).Value = 45008
$ws.Cells.Item(336, 5).Value = 13
$ws.Cells.Item(336, 6).Value = 100112001
$ws.Cells.Item(336, 7).Value = "Berenjena"
$ws.Cells.Item(336, 8).Value = "Sin especificar"
$ws.Cells.Item(336, 9).Value = "Primera"
$ws.Cells.Item(336, 10).Value = 90
$ws.Cells.Item(336, 11).Value = 7000
$ws.Cells.Item(336, 12).Value = 8000
$ws.Cells.Item(336, 13).Value = 7500
$ws.Cells.Item(336, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(336, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(336, 16).Value = 150
$ws.Cells.Item(336, 17).Value = 50
$ws.Cells.Item(336, 18).Value = "Hortaliza"
